$d = $word.ActiveDocument

# --- Edit 1: merge " income per round" + "." runs into a single run ---
$d.Content.Find.Execute(" income per round.", $false, $false, $false, $false, $false, $true, 1, $false, " income per round.", 2) | Out-Null

# --- Edit 2: remove the stray _GoBack bookmark after "per gold unit." ---
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

# --- Edit 3: expand the "Programming Notes" section at the end of the document ---
$lastPara = $d.Paragraphs.Last
$rng = $lastPara.Range

$xml = @'
<w:body xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:p><w:pPr><w:rPr><w:b/><w:u w:val="single"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:ind w:left="270" w:hanging="270"/></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">Programming Notes: </w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:ind w:left="540" w:hanging="270"/></w:pPr><w:r><w:t>&#8220;Games&#8221; button on main screen opens &#8220;Game Settings&#8221; and starts game setup mode.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr><w:ind w:left="900"/></w:pPr><w:r><w:t xml:space="preserve">Game Setup Mode returns </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Atmega</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> to normal startup settings except now in game setup mode.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr><w:ind w:left="900"/></w:pPr><w:r><w:t>Allows first person to submit changes to settings to lock in game settings.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr><w:ind w:left="900"/></w:pPr><w:r><w:t>If already in game setup mode, game settings are displayed and available choices for player characters are shown as tombstones.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr><w:ind w:left="900"/></w:pPr><w:r><w:t>Selecting a tomb</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t>stone locks player character and opens up a rename box that shows normal character name as default.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr><w:ind w:left="900"/></w:pPr><w:r><w:t xml:space="preserve">Selecting &#8220;Ready&#8221; takes player to main game display and adds to </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>players</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> ready count.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr><w:ind w:left="900"/></w:pPr><w:r><w:t>When player ready count = # of players in game settings, Game enters running mode.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:ind w:left="540" w:hanging="270"/></w:pPr><w:r><w:t xml:space="preserve">&#8220;Game Settings allows anyone to change settings until first person </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>submits</w:t></w:r><w:r><w:t>No</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> limit on rounds.  Play continues until only one is &#8220;</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Stayin</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> Alive.&#8221;</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:ind w:left="540" w:hanging="270"/></w:pPr></w:p><w:p/><w:p><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p></w:body>
'@

$rng.InsertXML($xml)
